$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert a new row at 397, shifting the existing rows 397-443 down to 398-444 ---
$ws.Rows.Item(397).Insert()

# The newly-inserted blank row 397 gets "default" formatting from the engine;
# restore the normal leave-card row formatting by copying formats from row 398
# (which is the former row 397, still carrying the standard style set).
$ws.Range("A398:K398").Copy()
$ws.Range("A397:K397").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# K396 and K397 should carry the date-formatted style (as used for the "as of"
# remark date column elsewhere in the sheet, e.g. K23) rather than the plain
# style.
$ws.Range("K23").Copy()
$ws.Range("K396:K397").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# Restore the calculated "EARNED " column formula on the new row (lost by the
# format-only paste above).
$ws.Cells.Item(397, 7).Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# --- Fill in the new leave-card entries ---
# Row 396: SP (Special Privilege) leave, 1.25 earned, remarks date 45132 (8/1/2023)
$ws.Cells.Item(396, 2).Value = "SP(1-0-0)"
$ws.Cells.Item(396, 3).Value = 1.25
$ws.Cells.Item(396, 11).Value = 45132

# Row 397: SL (Sick Leave), 1 day absence w/ pay, remarks date 45133 (8/2/2023)
$ws.Cells.Item(397, 2).Value = "SL(1-0-0)"
$ws.Cells.Item(397, 8).Value = 1
$ws.Cells.Item(397, 11).Value = 45133

# --- Resize Table1 to cover the new last row (444) ---
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A8:K444"))

# Re-assert the calculated-column formula text on the final row (the resize
# above rewrites it to a structured "[@EARNED]" reference); keep the original
# long-form structured reference used throughout the rest of the column.
$ws.Cells.Item(444, 7).Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# --- Update the selection to match the new cursor position ---
$ws.Range("I397").Select()
